# Reorders the existing 2011-2021 monthly rows so that each year's
# Oct/Nov/Dec figures lead that year's block (matches the refreshed
# upstream export), then appends the newly published 2022 (Oct-Dec-led)
# and 2023 Jan-Jul rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Snapshot the existing data rows (A2:D133) keyed by period ----
$rngExisting = $ws.Range("A2:D133")
$existing = $rngExisting.Value()

$map = @{}
for ($i = 1; $i -le 132; $i++) {
    $period = $existing[$i,1]
    $map[$period] = @($existing[$i,2], $existing[$i,3], $existing[$i,4])
}

# ---- 2. Build the rotated 2011-2021 order: Oct,Nov,Dec first, then Jan..Sep ----
$newOrder = @()
foreach ($y in 2011..2021) {
    $yearBlock = @()
    for ($m = 1; $m -le 12; $m++) {
        $yearBlock += ("{0}-{1:D2}" -f $y, $m)
    }
    $rotated = $yearBlock[9..11] + $yearBlock[0..8]
    $newOrder += $rotated
}

# ---- 3. Newly released periods: 2022 (Oct-Dec led) then 2023 Jan-Jul ----
$newMap = @{}
$newMap["2022-01"] = @(99.8, 99.8, 100)
$newMap["2022-02"] = @(100.5, 100.7, 100.1)
$newMap["2022-03"] = @(101.1, 101.4, 100.2)
$newMap["2022-04"] = @(100.6, 100.8, 100.2)
$newMap["2022-05"] = @(100.1, 100.1, 100.3)
$newMap["2022-06"] = @(100, 99.90000000000001, 100.3)
$newMap["2022-07"] = @(98.7, 98.3, 100.2)
$newMap["2022-08"] = @(98.8, 98.40000000000001, 99.90000000000001)
$newMap["2022-09"] = @(99.90000000000001, 99.8, 100.1)
$newMap["2022-10"] = @(100.2, 100.1, 100.5)
$newMap["2022-11"] = @(100.1, 100, 100.1)
$newMap["2022-12"] = @(99.5, 99.40000000000001, 99.8)
$newMap["2023-01"] = @(99.59999999999999, 99.5, 99.7)
$newMap["2023-02"] = @(100, 100.1, 99.7)
$newMap["2023-03"] = @(100, 100, 100)
$newMap["2023-04"] = @(99.5, 99.40000000000001, 99.7)
$newMap["2023-05"] = @(99.09999999999999, 98.8, 99.8)
$newMap["2023-06"] = @(99.2, 98.90000000000001, 99.8)
$newMap["2023-07"] = @(99.8, 99.59999999999999, 100.3)

foreach ($m in 10..12) { $newOrder += ("2022-{0:D2}" -f $m) }
foreach ($m in 1..9)   { $newOrder += ("2022-{0:D2}" -f $m) }
foreach ($m in 1..7)   { $newOrder += ("2023-{0:D2}" -f $m) }

foreach ($key in $newMap.Keys) { $map[$key] = $newMap[$key] }

# ---- 4. Write the reordered + extended data back out ----
$rowCount = $newOrder.Count()
$out = New-Object 'object[,]' $rowCount,4
for ($i = 0; $i -lt $rowCount; $i++) {
    $period = $newOrder[$i]
    $vals = $map[$period]
    $out[$i,0] = $period
    $out[$i,1] = $vals[0]
    $out[$i,2] = $vals[1]
    $out[$i,3] = $vals[2]
}

$lastRow = $rowCount + 1
$targetRange = $ws.Range("A2:D$lastRow")
$targetRange.Value = $out

# ---- 5. Make sure the new period (column A) cells carry the same style ----
#        ("bold, centred, boxed") as the rest of the date column.
$ws.Range("A2").Copy()
$ws.Range("A134:A$lastRow").PasteSpecial(-4122)
